# "Version 1p with colors"
#
# 1. Bump the fixed "last saved" date placeholder text from 30.08.2021 to
#    31.08.2021 everywhere it appears (the Slide Master and all 11 Custom
#    Layouts each carry their own "Datumsplatzhalter" placeholder shape).
# 2. Re-position the full-bleed picture ("Grafik 6") on slide 4.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Datumsplatzhalter*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "31.08.2021"

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

# Every Custom Layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# Slide 4, shape 3 ("Grafik 6") - move the picture
$slide4 = $p.Slides.Item(4)
$pic = $slide4.Shapes.Item(3)
$pic.Left = 501.3344881889764
$pic.Top = -131.02102362204724
